# Daily attendance processing - 2026-01-04 08:39:14
#
# Normalises the "Recorded By" column (G) on the "Session Analysis Results"
# sheet: wherever the literal recorder name "System" (exact case) appears
# together with other recorder names in the comma-separated list, "System"
# is moved to the front of the list while the remaining names keep their
# original relative order. Rows where "System" is the only recorder, or
# where "System" is not present at all, are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 through row 157 (row 1 is the header).
$firstRow = 2
$lastRow = 157

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "
    if ($parts.Count -le 1) { continue }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }
    if (-not $hasSystem) { continue }

    $others = $parts | Where-Object { -not $_.Equals("System") }
    $newParts = @("System") + $others
    $newVal = $newParts -join ", "

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
